$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15
$ws.Range("A15").Value = "No Emulator Data"
$ws.Range("C15").Value = "14"
$ws.Range("D15").Value = "3.0"
$ws.Range("E15").Value = "10.50"
$ws.Range("F15").Value = "26413414"
$ws.Range("G15").Value = "PayNow"
$ws.Range("H15").Value = "en_US"
$ws.Range("W15").Value = "udf data 4"
$ws.Range("X15").Value = "udf data 5"
$ws.Range("AB15").Value = "udf data 9"
$ws.Range("AC15").Value = "udf data 10"

# S column keeps the "Email" style (s=2) even when blank, matching the
# rest of the sheet -- copy the formatting from an existing S cell.
$ws.Range("S2").Copy()
$ws.Range("S15").PasteSpecial(-4122)

# Row 16
$ws.Range("A16").Value = "PayNow All Fields"
$ws.Range("C16").Value = "15"
$ws.Range("D16").Value = "3.0"
$ws.Range("E16").Value = "10.50"
$ws.Range("F16").Value = "26413255"
$ws.Range("G16").Value = "PayNow"
$ws.Range("H16").Value = "en_US"
$ws.Range("I16").Value = "Jasmine"
$ws.Range("J16").Value = "Patrinol"
$ws.Range("K16").Value = "258 Underwood rd"
$ws.Range("L16").Value = "Suite 600"
$ws.Range("M16").Value = "840"
$ws.Range("N16").Value = "Arlington"
$ws.Range("O16").Value = "VA"
$ws.Range("P16").Value = "22201"
$ws.Range("R16").Value = "Some Company"
$ws.Range("S2").Copy()
$ws.Range("S16").PasteSpecial(-4122)
$ws.Range("S16").Value = "iahmed@govolution.com"
$ws.Range("T16").Value = "udf data 1"
$ws.Range("U16").Value = "udf data 2"
$ws.Range("V16").Value = "udf data 3"
$ws.Range("W16").Value = "udf data 4"
$ws.Range("X16").Value = "udf data 5"
$ws.Range("Y16").Value = "udf data 6"
$ws.Range("Z16").Value = "Orange"
$ws.Range("AA16").Value = "Soccer"
$ws.Range("AB16").Value = "udf data 9"
$ws.Range("AC16").Value = "udf data 10"

# Row 17
$ws.Range("A17").Value = "No UnderPay"
$ws.Range("C17").Value = "16"
$ws.Range("D17").Value = "3.0"
$ws.Range("E17").Value = "10.50"
$ws.Range("F17").Value = "26413414"
$ws.Range("G17").Value = "PayNow"
$ws.Range("H17").Value = "en_US"
$ws.Range("I17").Value = "Curlonty"
$ws.Range("J17").Value = "Lachuga"
$ws.Range("K17").Value = "258 Underwood rd"
$ws.Range("L17").Value = "Suite 600"
$ws.Range("M17").Value = "840"
$ws.Range("N17").Value = "Arlington"
$ws.Range("O17").Value = "VA"
$ws.Range("P17").Value = "22201"
$ws.Range("S2").Copy()
$ws.Range("S17").PasteSpecial(-4122)
$ws.Range("T17").Value = "udf data 1"
$ws.Range("U17").Value = "udf data 2"

# Move the selection to the last-entered cell (matches a user who just
# typed the new row of data); the sheet scrolls back into view too.
$ws.Range("D17").Select()
